$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64 (ALC)
$ws.Range("H64").Value = 4000
$ws.Range("I64").Value = 4000
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 4000
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -3752
$ws.Range("N64").ClearContents()

# Row 67 (ALC)
$ws.Range("H67").Value = 4000
$ws.Range("I67").Value = 4000
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 4000
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -3142
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 23 (ARM)
$ws.Range("H23").Value = 28000.857

# Row 97 (ARM)
$ws.Range("H97").Value = 550.7895
$ws.Range("I97").Value = 471
$ws.Range("J97").Value = 850
$ws.Range("K97").Value = 471
$ws.Range("L97").Value = 850
$ws.Range("M97").Value = 25
$ws.Range("N97").Value = -1842

$ws = $wb.Worksheets.Item("BSM")
# Row 82 (BSM)
$ws.Range("H82").Value = 11051.417
$ws.Range("I82").Value = 4072.3333
$ws.Range("J82").Value = 31988.666
$ws.Range("K82").Value = 4072.3333
$ws.Range("L82").Value = 31988.666
$ws.Range("M82").Value = -3689.3333
$ws.Range("N82").Value = -32754.666

# Row 85 (BSM)
$ws.Range("H85").Value = 11051.417
$ws.Range("I85").Value = 4072.3333
$ws.Range("J85").Value = 31988.666
$ws.Range("K85").Value = 4072.3333
$ws.Range("L85").Value = 31988.666
$ws.Range("M85").Value = -2746.3333
$ws.Range("N85").Value = -34640.666

# Row 94 (BSM)
$ws.Range("H94").Value = 1616.2972
$ws.Range("I94").Value = 1483.8966
$ws.Range("J94").Value = 2096.25
$ws.Range("K94").Value = 1483.8966
$ws.Range("L94").Value = 2096.25
$ws.Range("M94").Value = -1032.8966
$ws.Range("N94").Value = -2998.25

# Row 97 (BSM)
$ws.Range("H97").Value = 13402.454
$ws.Range("I97").Value = 5485.4
$ws.Range("J97").Value = 20000
$ws.Range("K97").Value = 5485.4
$ws.Range("L97").Value = 20000
$ws.Range("M97").Value = -4494.4
$ws.Range("N97").Value = -21982

# Row 105 (BSM)
$ws.Range("H105").Value = 1841.4706
$ws.Range("I105").Value = 1656.6666
$ws.Range("J105").Value = 2285
$ws.Range("K105").Value = 1656.6666
$ws.Range("L105").Value = 2285
$ws.Range("M105").Value = 90.33339999999998
$ws.Range("N105").Value = -5779

# Row 107 (BSM)
$ws.Range("H107").Value = 5666.6665
$ws.Range("I107").Value = 4000
$ws.Range("J107").Value = 7333.3335
$ws.Range("K107").Value = 4000
$ws.Range("L107").Value = 7333.3335
$ws.Range("M107").Value = -2080
$ws.Range("N107").Value = -11173.3335

$ws = $wb.Worksheets.Item("CRP")
# Row 132 (CRP)
$ws.Range("H132").Value = 2584.7368
$ws.Range("I132").Value = 1515.0714
$ws.Range("J132").Value = 5579.8
$ws.Range("K132").Value = 4545.2142
$ws.Range("L132").Value = 16739.4
$ws.Range("M132").Value = -2015.2142
$ws.Range("N132").Value = -21799.4

$ws = $wb.Worksheets.Item("CUL")
# Row 34 (CUL)
$ws.Range("H34").Value = 5381.067
$ws.Range("J34").Value = 5752.2856
$ws.Range("L34").Value = 17256.8568
$ws.Range("N34").Value = -17424.8568

# Row 39 (CUL)
$ws.Range("H39").Value = 2649.6667
$ws.Range("I39").Value = 999
$ws.Range("J39").Value = 3475
$ws.Range("K39").Value = 2997
$ws.Range("L39").Value = 10425
$ws.Range("M39").Value = -2703
$ws.Range("N39").Value = -11013

# Row 55 (CUL)
$ws.Range("H55").Value = 2483.75
$ws.Range("J55").Value = 2900
$ws.Range("L55").Value = 8700
$ws.Range("N55").Value = -9054

# Row 64 (CUL)
$ws.Range("H64").Value = 9484.308000000001
$ws.Range("J64").Value = 11959.6
$ws.Range("L64").Value = 35878.8
$ws.Range("N64").Value = -36418.8

# Row 67 (CUL)
$ws.Range("H67").Value = 9484.308000000001
$ws.Range("J67").Value = 11959.6
$ws.Range("L67").Value = 35878.8
$ws.Range("N67").Value = -37750.8

# Row 69 (CUL)
$ws.Range("H69").Value = 16099.857
$ws.Range("I69").Value = 700
$ws.Range("K69").Value = 2100
$ws.Range("M69").Value = -1289

# Row 72 (CUL)
$ws.Range("H72").Value = 16099.857
$ws.Range("I72").Value = 700
$ws.Range("K72").Value = 6300
$ws.Range("M72").Value = -2244

# Row 113 (CUL)
$ws.Range("H113").Value = 1341.5
$ws.Range("I113").Value = 699.5
$ws.Range("J113").Value = 1469.9
$ws.Range("K113").Value = 2098.5
$ws.Range("L113").Value = 4409.700000000001
$ws.Range("M113").Value = 71.5
$ws.Range("N113").Value = -8749.700000000001

$ws = $wb.Worksheets.Item("GSM")
# Row 33 (GSM)
$ws.Range("H33").Value = 24178.334
$ws.Range("I33").Value = 517
$ws.Range("J33").Value = 36009
$ws.Range("K33").Value = 517
$ws.Range("L33").Value = 36009
$ws.Range("M33").Value = -265
$ws.Range("N33").Value = -36513

# Row 36 (GSM)
$ws.Range("H36").Value = 1339
$ws.Range("I36").Value = 1017
$ws.Range("J36").Value = 1500
$ws.Range("K36").Value = 1017
$ws.Range("L36").Value = 1500
$ws.Range("M36").Value = -532
$ws.Range("N36").Value = -2470

# Row 70 (GSM)
$ws.Range("H70").Value = 3860
$ws.Range("I70").Value = 3800
$ws.Range("J70").Value = 3980
$ws.Range("K70").Value = 3800
$ws.Range("L70").Value = 3980
$ws.Range("M70").Value = -3530
$ws.Range("N70").Value = -4520

# Row 73 (GSM)
$ws.Range("H73").Value = 3860
$ws.Range("I73").Value = 3800
$ws.Range("J73").Value = 3980
$ws.Range("K73").Value = 3800
$ws.Range("L73").Value = 3980
$ws.Range("M73").Value = -2864
$ws.Range("N73").Value = -5852

# Row 132 (GSM)
$ws.Range("H132").Value = 25644314
$ws.Range("I132").Value = 45456330
$ws.Range("J132").Value = 5235.1177
$ws.Range("K132").Value = 136368990
$ws.Range("L132").Value = 15705.3531
$ws.Range("M132").Value = -136366460
$ws.Range("N132").Value = -20765.3531

$ws = $wb.Worksheets.Item("LTW")
# Row 14 (LTW)
$ws.Range("H14").Value = 2985.3333
$ws.Range("I14").Value = 2985.3333
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 2985.3333
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -2813.3333
$ws.Range("N14").ClearContents()

# Row 132 (LTW)
$ws.Range("H132").Value = 3472.5
$ws.Range("I132").Value = 2181.5454
$ws.Range("J132").Value = 4763.4546
$ws.Range("K132").Value = 6544.6362
$ws.Range("L132").Value = 14290.3638
$ws.Range("M132").Value = -4014.6362
$ws.Range("N132").Value = -19350.3638

# Row 136 (LTW)
$ws.Range("H136").Value = 2264.1365
$ws.Range("I136").Value = 1765.3529
$ws.Range("J136").Value = 3960
$ws.Range("K136").Value = 5296.0587
$ws.Range("L136").Value = 11880
$ws.Range("M136").Value = -2746.0587
$ws.Range("N136").Value = -16980

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (WVR)
$ws.Range("H81").Value = 960
$ws.Range("I81").Value = 933.3333
$ws.Range("J81").Value = 971.4286
$ws.Range("K81").Value = 1866.6666
$ws.Range("L81").Value = 1942.8572
$ws.Range("M81").Value = -805.6666
$ws.Range("N81").Value = -4064.8572

# Row 84 (WVR)
$ws.Range("H84").Value = 960
$ws.Range("I84").Value = 933.3333
$ws.Range("J84").Value = 971.4286
$ws.Range("K84").Value = 9333.333000000001
$ws.Range("L84").Value = 9714.286
$ws.Range("M84").Value = -4029.333000000001
$ws.Range("N84").Value = -20322.286

# Row 132 (WVR)
$ws.Range("H132").Value = 6465.788
$ws.Range("I132").Value = 2759.4783
$ws.Range("J132").Value = 14990.3
$ws.Range("K132").Value = 8278.4349
$ws.Range("L132").Value = 44970.89999999999
$ws.Range("M132").Value = -5748.4349
$ws.Range("N132").Value = -50030.89999999999
